# Inserts a new weekly price record as row 308 on the "Hortaliza, Femacal de
# La Calera - Berenjena" sheet, pushing the former rows 308-333 down to
# 309-334 (the sheet's used range grows from A1:R333 to A1:R334).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a blank row at 308.
$ws.Rows.Item(308).Insert()

# Populate the new row with this week's record (values constant across the
# series carried over from the surrounding rows; the new observation's own
# figures come from the source data for the latest week).
$ws.Cells.Item(308, 1).Value = 3
$ws.Cells.Item(308, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(308, 3).Value = "Coquimbo"
$ws.Cells.Item(308, 4).Value = 44783
$ws.Cells.Item(308, 5).Value = 5
$ws.Cells.Item(308, 6).Value = 100112001
$ws.Cells.Item(308, 7).Value = "Berenjena"
$ws.Cells.Item(308, 8).Value = "Sin especificar"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 105
$ws.Cells.Item(308, 11).Value = 8500
$ws.Cells.Item(308, 12).Value = 9000
$ws.Cells.Item(308, 13).Value = 8762
$ws.Cells.Item(308, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(308, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(308, 16).Value = 146
$ws.Cells.Item(308, 17).Value = 60
$ws.Cells.Item(308, 18).Value = "Hortaliza"
